$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows above the current row 19 ("otherprefix:" row), shifting
# that row and the row below it down to rows 21 and 22.
$ws.Rows("19:20").Insert()

# Update the existing "very very deep path" value cell (J17) to its new text.
$ws.Range("J17").Value = "veryverydeeppath-value"

# Fill in the two newly inserted rows with the new overlapping JSON paths.
$ws.Range("A18").Value = "v.e.r.y.v.e.r.y.d.e.e.p"
$ws.Range("J18").Value = "veryverydeep-value"

$ws.Range("A19").Value = "v.e.r.y.v.e.r.y.d.e.e.p.v1"
$ws.Range("J19").Value = "veryverydeepv1-value"

# Register new defined names for the new JSON paths.
$wb.Names.Add("json.v.e.r.y.v.e.r.y.d.e.e.p", "=Sheet1!`$J`$18")
$wb.Names.Add("json.v.e.r.y.v.e.r.y.d.e.e.p.v1", "=Sheet1!`$J`$19")

# Leave the active selection on J19, matching the authored edit.
$ws.Range("J19").Select()
